$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 = "REX_DEF", matching the formatting of the existing
# header cells (e.g. E1 = "REX_DESC") by copying E1's format onto F1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "REX_DEF"
